$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "HINT" column (C) with its header and per-row hint values.
# (Write C3 before C2 so the shared-string table keeps the same
# insertion order as the authored workbook: HINT, c****er, s****er.)
$ws.Range("C1").Value = "HINT"
$ws.Range("C3").Value = "c****er"
$ws.Range("C2").Value = "s****er"

# Material-design-ish touch: widen column B so the Japanese prompts aren't clipped.
$ws.Columns("B").ColumnWidth = 36.83

# Move the active selection (as last left by the author) to D10.
$ws.Range("D10").Select() | Out-Null
